$wb = $excel.ActiveWorkbook

# --- Sheet "addVendor": shift AT-66..AT-70 / Auto_Vendor 66..70 -> AT-86..AT-90 / Auto_Vendor 86..90
$wsAdd = $wb.Worksheets.Item("addVendor")
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $new = 86 + $i
    $wsAdd.Range("D$row").Value = "AT-$new"
    $wsAdd.Range("E$row").Value = "Auto_Vendor $new"
}

# --- Sheet "editVendor": shift AT_EDT-71..75 / Auto_Vendor_edit_71..75 -> AT_EDT-91..95 / Auto_Vendor_edit_91..95
$wsEdit = $wb.Worksheets.Item("editVendor")
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $new = 91 + $i
    $wsEdit.Range("D$row").Value = "AT_EDT-$new"
    $wsEdit.Range("E$row").Value = "Auto_Vendor_edit_$new"
}

# --- Sheet "deleteVendor": shift AT_DEL_96..100 / Delete_Vendor_96..100 -> AT_DEL_116..120 / Delete_Vendor_116..120
$wsDel = $wb.Worksheets.Item("deleteVendor")
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $new = 116 + $i
    $wsDel.Range("D$row").Value = "AT_DEL_$new"
    $wsDel.Range("E$row").Value = "Delete_Vendor_$new"
}

# --- Sheet "syncVendor": shift TE-VE-IN-60..64 / Del_Vendor 60..64 -> TE-VE-IN-80..84 / Del_Vendor 80..84
$wsSync = $wb.Worksheets.Item("syncVendor")
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $new = 80 + $i
    $wsSync.Range("D$row").Value = "TE-VE-IN-$new"
    $wsSync.Range("E$row").Value = "Del_Vendor $new"
}

# --- Column width changes ---
# addVendor: set column C width (~10.78 chars wide)
$wsAdd.Columns.Item(3).ColumnWidth = 10

# deleteVendor: widen column D (~11.22 chars wide)
$wsDel.Columns.Item(4).ColumnWidth = 10.333333333333334

# --- syncVendor sheetView changes: clear horizontal scroll (topLeftCell), change selection ---
$wsSync.Activate()
$wsSync.Range("D2:E6").Select()
